$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the w:proofErr spellStart/spellEnd bracketing "KarHan".
#    The object model has no direct "remove proofing error marker"
#    call, so rebuild the "Team: KarHan" paragraph from scratch
#    (new paragraph with the same text/formatting) and delete the
#    old one, which drops the orphaned proofErr markers entirely.
# ------------------------------------------------------------------
$oldTeamPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "KarHan") {
        $oldTeamPara = $p
        break
    }
}

$insPt = $d.Range($oldTeamPara.Range.Start, $oldTeamPara.Range.Start)
$insPt.InsertParagraphBefore()

$newTeamPara = $d.Paragraphs(2)
$teamStart = $newTeamPara.Range.Start
$r1 = $d.Range($teamStart, $teamStart)
$r1.InsertAfter("Team: ")
$r2 = $d.Range($teamStart + 6, $teamStart + 6)
$r2.InsertAfter("KarHan")

$oldTeamParaNow = $d.Paragraphs(3)
$oldFull = $d.Range($oldTeamParaNow.Range.Start, $oldTeamParaNow.Range.End)
$oldFull.Delete()

# ------------------------------------------------------------------
# 2) Add a new acceptance paragraph "I accept this WBA – Kar Kei"
#    right after "I accept this WBA – Wei Han", keeping the existing
#    trailing empty paragraph intact, and mark the end of the new
#    text with a collapsed "_GoBack" bookmark (as Word stamps at the
#    point of the most recent edit).
# ------------------------------------------------------------------
$weiHanPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "I accept this WBA") {
        $weiHanPara = $p
        break
    }
}

$weiHanPara.Range.InsertParagraphAfter()

$karKeiStart = $weiHanPara.Range.End
$insKarKei = $d.Range($karKeiStart, $karKeiStart)
# Append a temporary trailing placeholder character; a collapsed
# range positioned exactly at a paragraph's text end (i.e. right
# before its paragraph mark) is mishandled by this host, so the
# bookmark is anchored one character early and the placeholder is
# trimmed off afterwards.
$insKarKei.InsertBefore("I accept this WBA – Kar KeiX")

$karKeiPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "Kar KeiX") {
        $karKeiPara = $p
        break
    }
}

$placeholderPos = $karKeiPara.Range.End - 2
$bmRange = $d.Range($placeholderPos, $placeholderPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$placeholderRange = $d.Range($placeholderPos, $placeholderPos + 1)
$placeholderRange.Text = ""

Write-Host "Done"
